# Adds a new survey wave ("26. 1. 2021") as the last column on both sheets,
# and bumps the "aktualizace" date in the two summary title cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": new column W (was last column V), rows 2-45 hold values
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Header cell: copy the format of the previous header (V1) onto W1, then set text
$wsData.Range("V1").Copy() | Out-Null
$wsData.Range("W1").PasteSpecial(-4122) | Out-Null
$wsData.Range("W1").Value = "26. 1. 2021"

$dataValues = @(0.23,0.12,0.5,0.3,0.15,0.23,0.34,0.18,0.24,0.22,0.3,0.36,0.19,0.21,0.26,0.21,0.33,0.27,0.17,0.15,0.14,0.21,0.45,0.43,0.1,0.07000000000000001,0.12,0.24,0.08,0.1,0.12,0.22,0.2,0.08,0.13,0.13,0.06,0.26,0.15,0.05,0.04,0.08,0.13,0.23)

for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 23).Value = $dataValues[$i]
}

# ---------------------------------------------------------------------
# Sheet "pocetR": new column V (was last column U), rows 2-23 hold values
# ---------------------------------------------------------------------
$wsCount = $wb.Worksheets.Item("pocetR")

$wsCount.Range("U1").Copy() | Out-Null
$wsCount.Range("V1").PasteSpecial(-4122) | Out-Null
$wsCount.Range("V1").Value = "26. 1. 2021"

$countValues = @(2131,237,491,1403,1013,190,621,307,966,176,132,857,984,726,421,272,829,675,203,526,385,242)

for ($i = 0; $i -lt $countValues.Length; $i++) {
    $row = $i + 2
    $wsCount.Cells.Item($row, 22).Value = $countValues[$i]
}

# Trailing blank cell on the title row, matching the other blank cells in that row.
# (A plain `.Value = ""` assignment is treated as "clear the cell" by the engine and
# would not persist an actual empty-string cell, so copy the existing blank sibling.)
$wsCount.Range("U24").Copy() | Out-Null
$wsCount.Range("V24").PasteSpecial(-4163) | Out-Null

# ---------------------------------------------------------------------
# Bump the "aktualizace" date in both summary title cells
# ---------------------------------------------------------------------
$wsData.Range("A46").Replace("12. 1. 2021", "2. 2. 2021") | Out-Null
$wsCount.Range("A24").Replace("12. 1. 2021", "2. 2. 2021") | Out-Null
